# Generate Report for Handoff
#
# The "34d0a952-eade-4acf-8ef8-b8771c6a7fdb.md" row (row 6 in each table)
# gets a refreshed "Handoff" generation timestamp across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-25 10:43:32"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-25 10:43:28"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-25 10:43:32"
